# Insert a new row of data at row 353 (shifting existing rows 353:451 down to 354:452)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 353; this pushes all rows
# 353-451 down to 354-452 and the sheet dimension grows to A1:R452.
$ws.Rows("353:353").Insert()

# Populate the newly inserted row 353 with the new record's values.
$ws.Range("A353").Value = 3
$ws.Range("B353").Value = "Femacal de La Calera"
$ws.Range("C353").Value = "Coquimbo"
$ws.Range("D353").Value = 44841
$ws.Range("E353").Value = 5
$ws.Range("F353").Value = 100112031
$ws.Range("G353").Value = "Poroto verde"
$ws.Range("H353").Value = "Magnum"
$ws.Range("I353").Value = "Primera"
$ws.Range("J353").Value = 78
$ws.Range("K353").Value = 33000
$ws.Range("L353").Value = 34000
$ws.Range("M353").Value = 33513
$ws.Range("N353").Value = "$/saco 25 kilos"
$ws.Range("O353").Value = "Región de Arica y Parinacota"
$ws.Range("P353").Value = 1341
$ws.Range("Q353").Value = 25
$ws.Range("R353").Value = "Hortaliza"
